$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.813602089881897
$ws.Range("B1").Value = 2.789194822311401
$ws.Range("C1").Value = 4.754696846008301
$ws.Range("D1").Value = 2.783067226409912
$ws.Range("E1").Value = 1.346285820007324
